$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 77030.16
$ws.Range("I33").Value = 83441
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 83441
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -83212
$ws.Range("H38").Value = 144
$ws.Range("I38").Value = 144
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 432
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -60
$ws.Range("N38").ClearContents()
$ws.Range("H76").Value = 23408182
$ws.Range("I76").Value = 39289460
$ws.Range("J76").Value = 4194.1055
$ws.Range("K76").Value = 39289460
$ws.Range("L76").Value = 4194.1055
$ws.Range("M76").Value = -39289145
$ws.Range("N76").Value = -4824.1055
$ws.Range("H79").Value = 23408182
$ws.Range("I79").Value = 39289460
$ws.Range("J79").Value = 4194.1055
$ws.Range("K79").Value = 39289460
$ws.Range("L79").Value = 4194.1055
$ws.Range("M79").Value = -39288368
$ws.Range("N79").Value = -6378.1055
$ws.Range("H129").Value = 1764570.4
$ws.Range("I129").Value = 300
$ws.Range("J129").Value = 2315905
$ws.Range("K129").Value = 900
$ws.Range("L129").Value = 6947715
$ws.Range("M129").Value = 4100
$ws.Range("N129").Value = -6957715
$ws.Range("H132").Value = 627967
$ws.Range("I132").Value = 3751.3635
$ws.Range("J132").Value = 2001241.4
$ws.Range("K132").Value = 11254.0905
$ws.Range("L132").Value = 6003724.199999999
$ws.Range("M132").Value = -8724.0905
$ws.Range("H136").Value = 30766.666
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 30766.666
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 30766.666
$ws.Range("N136").Value = -40966.666
$ws.Range("H137").Value = 3428.8572
$ws.Range("I137").Value = 3083.6667
$ws.Range("J137").Value = 5500
$ws.Range("K137").Value = 9251.000100000001
$ws.Range("L137").Value = 16500
$ws.Range("M137").Value = -6701.000100000001
$ws.Range("N137").Value = -21600
$ws.Range("H138").Value = 4499.0806
$ws.Range("I138").Value = 1037.8334
$ws.Range("J138").Value = 4976.494
$ws.Range("K138").Value = 3113.5002
$ws.Range("L138").Value = 14929.482
$ws.Range("M138").Value = 2026.4998
$ws.Range("N138").Value = -25209.482
$ws.Range("H139").Value = 34132
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 34132
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 34132
$ws.Range("N139").Value = -44412
$ws.Range("H140").Value = 47737.625
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 47737.625
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 47737.625
$ws.Range("N140").Value = -58097.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 8500
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 8500
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 8500
$ws.Range("N9").Value = -8840
$ws.Range("H20").Value = 8500
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 8500
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 8500
$ws.Range("N20").Value = -9040
$ws.Range("H32").Value = 7408.6665
$ws.Range("I32").Value = 7066.9
$ws.Range("J32").Value = 10142.8
$ws.Range("K32").Value = 7066.9
$ws.Range("L32").Value = 10142.8
$ws.Range("M32").Value = -6779.9
$ws.Range("N32").Value = -10716.8
$ws.Range("H33").Value = 3514.5
$ws.Range("I33").Value = 3000
$ws.Range("J33").Value = 4029
$ws.Range("K33").Value = 3000
$ws.Range("L33").Value = 4029
$ws.Range("M33").Value = -2671
$ws.Range("N33").Value = -4687
$ws.Range("H45").Value = 2580
$ws.Range("I45").Value = 2176.6667
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 2176.6667
$ws.Range("L45").Value = 5000
$ws.Range("M45").Value = -1799.6667
$ws.Range("N45").Value = -5754
$ws.Range("H102").Value = 5266.6665
$ws.Range("I102").Value = 6400
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 6400
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -4778
$ws.Range("N102").Value = -6244
$ws.Range("H132").Value = 2085.639
$ws.Range("I132").Value = 1511.2333
$ws.Range("J132").Value = 4957.6665
$ws.Range("K132").Value = 4533.699900000001
$ws.Range("L132").Value = 14872.9995
$ws.Range("M132").Value = -2003.699900000001
$ws.Range("H139").Value = 24500
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 24500
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 24500
$ws.Range("N139").Value = -34780

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1219.4
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 1219.4
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 1219.4
$ws.Range("N64").Value = -1669.4
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 1219.4
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 1219.4
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 1219.4
$ws.Range("N67").Value = -2779.4
$ws.Range("M67").ClearContents()
$ws.Range("H105").Value = 3701.2727
$ws.Range("I105").Value = 2871.4
$ws.Range("J105").Value = 12000
$ws.Range("K105").Value = 2871.4
$ws.Range("L105").Value = 12000
$ws.Range("M105").Value = -1124.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3176.5
$ws.Range("I31").Value = 1177.0454
$ws.Range("J31").Value = 8675
$ws.Range("K31").Value = 1177.0454
$ws.Range("L31").Value = 8675
$ws.Range("M31").Value = -882.0454
$ws.Range("N31").Value = -9265
$ws.Range("H34").Value = 3176.5
$ws.Range("I34").Value = 1177.0454
$ws.Range("J34").Value = 8675
$ws.Range("K34").Value = 1177.0454
$ws.Range("L34").Value = 8675
$ws.Range("M34").Value = -975.0454
$ws.Range("N34").Value = -9079
$ws.Range("H86").Value = 3952.625
$ws.Range("I86").Value = 4024.2
$ws.Range("J86").Value = 3833.3333
$ws.Range("K86").Value = 4024.2
$ws.Range("L86").Value = 3833.3333
$ws.Range("M86").Value = -2901.2
$ws.Range("N86").Value = -6079.3333
$ws.Range("H89").Value = 3952.625
$ws.Range("I89").Value = 4024.2
$ws.Range("J89").Value = 3833.3333
$ws.Range("K89").Value = 20121
$ws.Range("L89").Value = 19166.6665
$ws.Range("M89").Value = -14505
$ws.Range("N89").Value = -30398.6665
$ws.Range("H132").Value = 2408.6086
$ws.Range("I132").Value = 1458.5714
$ws.Range("J132").Value = 2824.25
$ws.Range("K132").Value = 4375.7142
$ws.Range("L132").Value = 8472.75
$ws.Range("M132").Value = -1845.7142
$ws.Range("H138").Value = 41145
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 41145
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 41145
$ws.Range("N138").Value = -51425
$ws.Range("H140").Value = 53640
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 53640
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 53640
$ws.Range("N140").Value = -64000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 97.44444
$ws.Range("I8").Value = 97.44444
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 292.33332
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -153.33332
$ws.Range("H97").Value = 14286362
$ws.Range("I97").Value = 23809822
$ws.Range("J97").Value = 1170
$ws.Range("K97").Value = 71429466
$ws.Range("L97").Value = 3510
$ws.Range("M97").Value = -71428970
$ws.Range("N97").Value = -4502
$ws.Range("H98").Value = 1672.4375
$ws.Range("I98").Value = 700
$ws.Range("J98").Value = 1996.5834
$ws.Range("K98").Value = 2100
$ws.Range("L98").Value = 5989.7502
$ws.Range("M98").Value = -602
$ws.Range("N98").Value = -8985.7502
$ws.Range("H107").Value = 486800.12
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 519186.8
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 1557560.4
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -1561400.4
$ws.Range("H131").Value = 40929.58
$ws.Range("I131").Value = 200366
$ws.Range("J131").Value = 2968.524
$ws.Range("K131").Value = 601098
$ws.Range("L131").Value = 8905.572
$ws.Range("M131").Value = -596058
$ws.Range("N131").Value = -18985.572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4339.4287
$ws.Range("I70").Value = 4077.9546
$ws.Range("J70").Value = 5298.1665
$ws.Range("K70").Value = 4077.9546
$ws.Range("L70").Value = 5298.1665
$ws.Range("M70").Value = -3807.9546
$ws.Range("N70").Value = -5838.1665
$ws.Range("H73").Value = 4339.4287
$ws.Range("I73").Value = 4077.9546
$ws.Range("J73").Value = 5298.1665
$ws.Range("K73").Value = 4077.9546
$ws.Range("L73").Value = 5298.1665
$ws.Range("M73").Value = -3141.9546
$ws.Range("N73").Value = -7170.1665
$ws.Range("H113").Value = 2144.7144
$ws.Range("I113").Value = 2333.3333
$ws.Range("J113").Value = 1013
$ws.Range("K113").Value = 2333.3333
$ws.Range("L113").Value = 1013
$ws.Range("M113").Value = -163.3332999999998
$ws.Range("N113").Value = -5353
$ws.Range("H136").Value = 22430.6
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 22430.6
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 67291.79999999999
$ws.Range("N136").Value = -72391.79999999999
$ws.Range("H138").Value = 30660
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 30660
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 30660
$ws.Range("N138").Value = -40940
$ws.Range("H139").Value = 27052
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 27052
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 27052
$ws.Range("N139").Value = -37332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10645.818
$ws.Range("I61").Value = 13888
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 13888
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -13686
$ws.Range("H113").Value = 10645.818
$ws.Range("I113").Value = 13888
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 13888
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -11718
$ws.Range("H135").Value = 50333.332
$ws.Range("I135").Value = 24000
$ws.Range("J135").Value = 53625
$ws.Range("K135").Value = 24000
$ws.Range("L135").Value = 53625
$ws.Range("M135").Value = -18930
$ws.Range("N135").Value = -63765
$ws.Range("H138").Value = 24957.25
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 24957.25
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 24957.25
$ws.Range("N138").Value = -35237.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1908.3334
$ws.Range("I81").Value = 1596.4286
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3192.8572
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2131.8572
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1908.3334
$ws.Range("I84").Value = 1596.4286
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 15964.286
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -10660.286
$ws.Range("N84").Value = -40608
$ws.Range("H138").Value = 39640
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 39640
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 39640
$ws.Range("N138").Value = -49920
